# Generate Report for Handback
# Updates the handback-status report with freshly generated timestamps for
# the c617b04f-... row (row 3) on each sheet: the Overview sheet's
# "Latest HO Xliff Generate Date" column, and the per-locale sheets'
# "Correspond Handoff Datetime" / "Correspond Handback DateTime" columns.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: refresh the Latest HO Xliff Generate Date for row 3 ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G3").Value = "2016-09-05 18:54:47"

# --- zh-cn sheet: refresh handoff/handback datetimes for row 3 ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H3").Value = "2016-09-05 18:54:41"
$zhcn.Range("K3").Value = "2016-09-05 18:55:13"

# --- de-de sheet: refresh handoff/handback datetimes for row 3 ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H3").Value = "2016-09-05 18:54:47"
$dede.Range("K3").Value = "2016-09-05 18:55:21"
